# Refresh market-board derived profit columns (H/I/J/K/L/M/N) on several
# Leve rows across the class worksheets, per the scheduled market-data run.
$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

# ALC row 15: Morning Glass of Ether
$ws_ALC.Range("H15").Value = 1324.6364
$ws_ALC.Range("I15").Value = 1324.6364
$ws_ALC.Range("K15").Value = 3973.9092
$ws_ALC.Range("M15").Value = -3804.9092

# ALC row 96: Scroll Down
$ws_ALC.Range("H96").Value = 2425
$ws_ALC.Range("I96").Value = 2425
$ws_ALC.Range("K96").Value = 7275
$ws_ALC.Range("M96").Value = -5902

# ALC row 103: Let Loose the Juice
$ws_ALC.Range("H103").Value = 639.75
$ws_ALC.Range("J103").Value = 639.75
$ws_ALC.Range("L103").Value = 1919.25
$ws_ALC.Range("N103").Value = -3091.25

# ALC row 112: Making Ends Meet
$ws_ALC.Range("H112").Value = 1828.6666
$ws_ALC.Range("J112").Value = 1985.8334
$ws_ALC.Range("L112").Value = 5957.5002
$ws_ALC.Range("N112").Value = -8173.5002

# ALC row 138: All-night Crafting
$ws_ALC.Range("H138").Value = 2602.9333
$ws_ALC.Range("I138").Value = 1577.7142
$ws_ALC.Range("K138").Value = 4733.142599999999
$ws_ALC.Range("M138").Value = 406.8574000000008

# ARM row 2: Ain't Got No Ingots
$ws_ARM.Range("H2").Value = 1344.25
$ws_ARM.Range("I2").Value = 769.75
$ws_ARM.Range("K2").Value = 769.75
$ws_ARM.Range("M2").Value = -656.75

# ARM row 32: Ingot We Trust
$ws_ARM.Range("H32").Value = 3778.7354
$ws_ARM.Range("I32").Value = 3778.7354
$ws_ARM.Range("J32").Value = 0
$ws_ARM.Range("K32").Value = 3778.7354
$ws_ARM.Range("L32").Value = 0
$ws_ARM.Range("M32").Value = -3491.7354
$ws_ARM.Range("N32").ClearContents()

# ARM row 97: Ore for Me
$ws_ARM.Range("H97").Value = 729
$ws_ARM.Range("I97").Value = 688.2143
$ws_ARM.Range("K97").Value = 688.2143
$ws_ARM.Range("M97").Value = -192.2143

# ARM row 102: Smells of Rich Tama-hagane
$ws_ARM.Range("H102").Value = 3275
$ws_ARM.Range("I102").Value = 1595
$ws_ARM.Range("J102").Value = 9995
$ws_ARM.Range("K102").Value = 1595
$ws_ARM.Range("L102").Value = 9995
$ws_ARM.Range("M102").Value = 27
$ws_ARM.Range("N102").Value = -13239

# ARM row 107: Shielding the Realm
$ws_ARM.Range("H107").Value = 20000
$ws_ARM.Range("J107").Value = 20000
$ws_ARM.Range("L107").Value = 20000
$ws_ARM.Range("N107").Value = -27680

# ARM row 110: Scheduled Maintenance
$ws_ARM.Range("H110").Value = 2566.889
$ws_ARM.Range("I110").Value = 2586.7144
$ws_ARM.Range("J110").Value = 2497.5
$ws_ARM.Range("K110").Value = 2586.7144
$ws_ARM.Range("L110").Value = 2497.5
$ws_ARM.Range("M110").Value = -541.7143999999998
$ws_ARM.Range("N110").Value = -6587.5

# ARM row 116: No Scope
$ws_ARM.Range("H116").Value = 1344.25
$ws_ARM.Range("I116").Value = 769.75
$ws_ARM.Range("K116").Value = 769.75
$ws_ARM.Range("M116").Value = 1524.25

# ARM row 132: Don't Bore Me, Ore Me
$ws_ARM.Range("H132").Value = 1694.5834
$ws_ARM.Range("J132").Value = 2724.75
$ws_ARM.Range("L132").Value = 8174.25
$ws_ARM.Range("N132").Value = -13234.25

# BSM row 3: Hells Bells
$ws_BSM.Range("H3").Value = 1344.25
$ws_BSM.Range("I3").Value = 769.75
$ws_BSM.Range("K3").Value = 769.75
$ws_BSM.Range("M3").Value = -655.75

# BSM row 94: High Steal
$ws_BSM.Range("H94").Value = 1721.1052
$ws_BSM.Range("I94").Value = 1688.2941
$ws_BSM.Range("K94").Value = 1688.2941
$ws_BSM.Range("M94").Value = -1237.2941

# BSM row 105: Ingot to Wing It
$ws_BSM.Range("H105").Value = 2000
$ws_BSM.Range("I105").Value = 1850
$ws_BSM.Range("J105").Value = 2150
$ws_BSM.Range("K105").Value = 1850
$ws_BSM.Range("L105").Value = 2150
$ws_BSM.Range("M105").Value = -103
$ws_BSM.Range("N105").Value = -5644

# CRP row 7: Gridania's Got Talent
$ws_CRP.Range("H7").Value = 60.833332
$ws_CRP.Range("I7").Value = 72.59999999999999
$ws_CRP.Range("J7").Value = 2
$ws_CRP.Range("K7").Value = 72.59999999999999
$ws_CRP.Range("L7").Value = 2
$ws_CRP.Range("M7").Value = 40.40000000000001
$ws_CRP.Range("N7").Value = -228

# CRP row 99: O Pine
$ws_CRP.Range("H99").Value = 7150.7
$ws_CRP.Range("I99").Value = 7150.7
$ws_CRP.Range("J99").Value = 0
$ws_CRP.Range("K99").Value = 7150.7
$ws_CRP.Range("L99").Value = 0
$ws_CRP.Range("M99").Value = -5652.7
$ws_CRP.Range("N99").ClearContents()

# CRP row 105: Zelkova, My Love
$ws_CRP.Range("H105").Value = 1627.8334
$ws_CRP.Range("J105").Value = 2746
$ws_CRP.Range("L105").Value = 2746
$ws_CRP.Range("N105").Value = -6240

# CRP row 126: A Better Conductor
$ws_CRP.Range("H126").Value = 7150.7
$ws_CRP.Range("I126").Value = 7150.7
$ws_CRP.Range("J126").Value = 0
$ws_CRP.Range("K126").Value = 21452.1
$ws_CRP.Range("L126").Value = 0
$ws_CRP.Range("M126").Value = -18982.1
$ws_CRP.Range("N126").ClearContents()

# CUL row 107: Slippery Service
$ws_CUL.Range("H107").Value = 468.16666
$ws_CUL.Range("I107").Value = 402
$ws_CUL.Range("K107").Value = 1206
$ws_CUL.Range("M107").Value = 714

# CUL row 121: A Cookie for Your Troubles
$ws_CUL.Range("H121").Value = 764.05554
$ws_CUL.Range("J121").Value = 981
$ws_CUL.Range("L121").Value = 2943
$ws_CUL.Range("N121").Value = -5563

# CUL row 132: More Mezcal
$ws_CUL.Range("H132").Value = 623
$ws_CUL.Range("I132").Value = 199
$ws_CUL.Range("J132").Value = 764.3333
$ws_CUL.Range("K132").Value = 1791
$ws_CUL.Range("L132").Value = 6878.9997
$ws_CUL.Range("M132").Value = 739
$ws_CUL.Range("N132").Value = -11938.9997

# GSM row 33: Thaumaturge Is Magic
$ws_GSM.Range("H33").Value = 19667
$ws_GSM.Range("J33").Value = 19667
$ws_GSM.Range("L33").Value = 19667
$ws_GSM.Range("N33").Value = -20171

# GSM row 59: Sew Not Doing This
$ws_GSM.Range("H59").Value = 0
$ws_GSM.Range("I59").Value = 0
$ws_GSM.Range("J59").Value = 0
$ws_GSM.Range("K59").Value = 0
$ws_GSM.Range("L59").Value = 0
$ws_GSM.Range("M59").ClearContents()
$ws_GSM.Range("N59").ClearContents()

# GSM row 132: On Board for Lar
$ws_GSM.Range("H132").Value = 3474
$ws_GSM.Range("I132").Value = 3474
$ws_GSM.Range("K132").Value = 10422
$ws_GSM.Range("M132").Value = -7892

# GSM row 135: Fan of the Foreign
$ws_GSM.Range("H135").Value = 48779.5
$ws_GSM.Range("J135").Value = 48779.5
$ws_GSM.Range("L135").Value = 48779.5
$ws_GSM.Range("N135").Value = -58919.5

# LTW row 22: Skin off Their Backs
$ws_LTW.Range("H22").Value = 1879.9
$ws_LTW.Range("I22").Value = 1879.9
$ws_LTW.Range("K22").Value = 1879.9
$ws_LTW.Range("M22").Value = -1584.9

# LTW row 27: Fire and Hide
$ws_LTW.Range("H27").Value = 1879.9
$ws_LTW.Range("I27").Value = 1879.9
$ws_LTW.Range("K27").Value = 1879.9
$ws_LTW.Range("M27").Value = -1772.9

# LTW row 47: Springtime for Coerthas
$ws_LTW.Range("H47").Value = 30031
$ws_LTW.Range("I47").Value = 0
$ws_LTW.Range("J47").Value = 30031
$ws_LTW.Range("K47").Value = 0
$ws_LTW.Range("L47").Value = 30031
$ws_LTW.Range("M47").ClearContents()
$ws_LTW.Range("N47").Value = -31011

# LTW row 52: The Tao of Rabbits
$ws_LTW.Range("H52").Value = 30031
$ws_LTW.Range("I52").Value = 0
$ws_LTW.Range("J52").Value = 30031
$ws_LTW.Range("K52").Value = 0
$ws_LTW.Range("L52").Value = 30031
$ws_LTW.Range("M52").ClearContents()
$ws_LTW.Range("N52").Value = -30497

# LTW row 68: You Could Say It's a Moving Target
$ws_LTW.Range("H68").Value = 1824.7142
$ws_LTW.Range("I68").Value = 1832.1666
$ws_LTW.Range("J68").Value = 1780
$ws_LTW.Range("K68").Value = 1832.1666
$ws_LTW.Range("L68").Value = 1780
$ws_LTW.Range("M68").Value = -1083.1666
$ws_LTW.Range("N68").Value = -3278

# LTW row 71: They Call It Bloody Mary (L)
$ws_LTW.Range("H71").Value = 1824.7142
$ws_LTW.Range("I71").Value = 1832.1666
$ws_LTW.Range("J71").Value = 1780
$ws_LTW.Range("K71").Value = 9160.833000000001
$ws_LTW.Range("L71").Value = 8900
$ws_LTW.Range("M71").Value = -5416.833000000001

# LTW row 100: Tiger in the Sack
$ws_LTW.Range("H100").Value = 1247.6666
$ws_LTW.Range("I100").Value = 1337.2
$ws_LTW.Range("J100").Value = 800
$ws_LTW.Range("K100").Value = 1337.2
$ws_LTW.Range("L100").Value = 800
$ws_LTW.Range("M100").Value = -796.2
$ws_LTW.Range("N100").Value = -1882

# LTW row 104: Brace Yourselves
$ws_LTW.Range("H104").Value = 22000
$ws_LTW.Range("J104").Value = 22000
$ws_LTW.Range("L104").Value = 22000
$ws_LTW.Range("N104").Value = -28988

# LTW row 136: Respect for Br'aax
$ws_LTW.Range("H136").Value = 2902.2
$ws_LTW.Range("I136").Value = 2902.2
$ws_LTW.Range("K136").Value = 8706.599999999999
$ws_LTW.Range("M136").Value = -6156.599999999999

# WVR row 62: Pride Up in Smoke
$ws_WVR.Range("H62").Value = 3870.2
$ws_WVR.Range("I62").Value = 3712.375
$ws_WVR.Range("K62").Value = 3712.375
$ws_WVR.Range("M62").Value = -3088.375

# WVR row 65: Desperate for Diversionaries (L)
$ws_WVR.Range("H65").Value = 3870.2
$ws_WVR.Range("I65").Value = 3712.375
$ws_WVR.Range("K65").Value = 18561.875
$ws_WVR.Range("M65").Value = -15441.875

# WVR row 107: Flax Wax
$ws_WVR.Range("H107").Value = 337
$ws_WVR.Range("I107").Value = 271.8
$ws_WVR.Range("K107").Value = 815.4000000000001
$ws_WVR.Range("M107").Value = 1104.6
